# Update of league bases, data edit for Poland Ekstraklasa sheet.
# Applies the per-cell deltas described by the upstream diff:
#  - fills in match-result (FTHG/FTAG/FTR) and recomputed odds columns for
#    rows 230-233 (matches that have since been played),
#  - refreshes the closing-odds columns for row 234,
#  - appends a brand-new fixture as row 235.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 230 ----
$ws.Range("H230").Value = 1
$ws.Range("I230").Value = 1
$ws.Range("J230").Value = "D"
$ws.Range("N230").Value = 2.9
$ws.Range("P230").Value = 2.4
$ws.Range("Q230").Value = 0.25
$ws.Range("R230").Value = 1.775
$ws.Range("S230").Value = 2.1
$ws.Range("U230").Value = 1.9
$ws.Range("V230").Value = 1.95
$ws.Range("W230").Value = -1
$ws.Range("X230").Value = 2.4
$ws.Range("Y230").Value = -1
$ws.Range("Z230").Value = 0.3875
$ws.Range("AA230").Value = -0.5
$ws.Range("AB230").Value = -0.5
$ws.Range("AC230").Value = 0.475

# ---- Row 231 ----
$ws.Range("H231").Value = 0
$ws.Range("I231").Value = 0
$ws.Range("J231").Value = "D"
$ws.Range("N231").Value = 4.75
$ws.Range("O231").Value = 3.5
$ws.Range("P231").Value = 1.8
$ws.Range("R231").Value = 2.025
$ws.Range("S231").Value = 1.825
$ws.Range("U231").Value = 1.85
$ws.Range("V231").Value = 2
$ws.Range("W231").Value = -1
$ws.Range("X231").Value = 2.5
$ws.Range("Y231").Value = -1
$ws.Range("Z231").Value = 1.025
$ws.Range("AA231").Value = -1
$ws.Range("AB231").Value = -1
$ws.Range("AC231").Value = 1

# ---- Row 232 ----
$ws.Range("H232").Value = 3
$ws.Range("I232").Value = 1
$ws.Range("J232").Value = "H"
$ws.Range("N232").Value = 2.1
$ws.Range("O232").Value = 3.3
$ws.Range("P232").Value = 3.6
$ws.Range("R232").Value = 1.8
$ws.Range("S232").Value = 2.05
$ws.Range("W232").Value = 1.1
$ws.Range("X232").Value = -1
$ws.Range("Y232").Value = -1
$ws.Range("Z232").Value = 0.8
$ws.Range("AA232").Value = -1
$ws.Range("AB232").Value = 1.025
$ws.Range("AC232").Value = -1

# ---- Row 233 ----
$ws.Range("H233").Value = 1
$ws.Range("I233").Value = 3
$ws.Range("J233").Value = "A"
$ws.Range("N233").Value = 3.5
$ws.Range("O233").Value = 3.4
$ws.Range("P233").Value = 2.15
$ws.Range("Q233").Value = 0.25
$ws.Range("R233").Value = 2
$ws.Range("S233").Value = 1.85
$ws.Range("U233").Value = 2.05
$ws.Range("V233").Value = 1.8
$ws.Range("W233").Value = -1
$ws.Range("X233").Value = -1
$ws.Range("Y233").Value = 1.15
$ws.Range("Z233").Value = -1
$ws.Range("AA233").Value = 0.8500000000000001
$ws.Range("AB233").Value = 1.05
$ws.Range("AC233").Value = -1

# ---- Row 234 ----
$ws.Range("N234").Value = 2.9
$ws.Range("O234").Value = 3
$ws.Range("P234").Value = 2.625
$ws.Range("T234").Value = 2
$ws.Range("U234").Value = 1.8
$ws.Range("V234").Value = 2.05

# ---- Row 235 (new fixture) ----
# Copy the number formatting used by the "id" and "Date" columns so the new
# row matches the look of the rest of the table.
$ws.Range("A234").Copy()
$ws.Range("A235").PasteSpecial(-4122)
$ws.Range("E234").Copy()
$ws.Range("E235").PasteSpecial(-4122)

$ws.Range("A235").Value = 233
$ws.Range("B235").Value = 6774470
$ws.Range("C235").Value = "Poland Ekstraklasa"
$ws.Range("D235").Value = "Poland Ekstraklasa"
$ws.Range("E235").Value = 45387.54166666666
$ws.Range("F235").Value = "Cracovia Krakow"
$ws.Range("G235").Value = "LKS Lodz"
$ws.Range("K235").Value = 1.5
$ws.Range("L235").Value = 4.333
$ws.Range("M235").Value = 6.5
$ws.Range("N235").Value = 1.5
$ws.Range("O235").Value = 4.333
$ws.Range("P235").Value = 6.5
$ws.Range("Q235").Value = -1
$ws.Range("R235").Value = 1.8
$ws.Range("S235").Value = 2.05
$ws.Range("T235").Value = 2.75
$ws.Range("U235").Value = 2.025
$ws.Range("V235").Value = 1.825
$ws.Range("W235").Value = 0
$ws.Range("X235").Value = 0
$ws.Range("Y235").Value = 0
$ws.Range("Z235").Value = 0
$ws.Range("AA235").Value = 0
